# Auto-generated edit script applying the Adamantoise_Profits diff
# Updates profit-calculation columns (H-N) across multiple job sheets
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H31").Value = 999
$ws_ALC.Range("J31").Value = 999
$ws_ALC.Range("L31").Value = 2997
$ws_ALC.Range("N31").Value = -3457
$ws_ALC.Range("H40").Value = 123689.6
$ws_ALC.Range("I40").Value = 1001663.3
$ws_ALC.Range("J40").Value = 3965.9092
$ws_ALC.Range("K40").Value = 1001663.3
$ws_ALC.Range("L40").Value = 3965.9092
$ws_ALC.Range("M40").Value = -1001488.3
$ws_ALC.Range("N40").Value = -4315.9092
$ws_ALC.Range("H43").Value = 7480.467
$ws_ALC.Range("J43").Value = 5632.4546
$ws_ALC.Range("L43").Value = 5632.4546
$ws_ALC.Range("N43").Value = -5770.4546
$ws_ALC.Range("H81").Value = 59324
$ws_ALC.Range("J81").Value = 59324
$ws_ALC.Range("L81").Value = 59324
$ws_ALC.Range("N81").Value = -61320
$ws_ALC.Range("H84").Value = 59324
$ws_ALC.Range("J84").Value = 59324
$ws_ALC.Range("L84").Value = 177972
$ws_ALC.Range("N84").Value = -187956
$ws_ALC.Range("H101").Value = 999
$ws_ALC.Range("I101").Value = 0
$ws_ALC.Range("J101").Value = 999
$ws_ALC.Range("K101").Value = 0
$ws_ALC.Range("L101").Value = 2997
$ws_ALC.Range("M101").ClearContents()
$ws_ALC.Range("N101").Value = -6241
$ws_ALC.Range("H125").Value = 7866.6665
$ws_ALC.Range("I125").Value = 1800
$ws_ALC.Range("K125").Value = 16200
$ws_ALC.Range("M125").Value = -13740
$ws_ALC.Range("H129").Value = 1136.1177
$ws_ALC.Range("I129").Value = 732.61536
$ws_ALC.Range("J129").Value = 2447.5
$ws_ALC.Range("K129").Value = 2197.84608
$ws_ALC.Range("L129").Value = 7342.5
$ws_ALC.Range("M129").Value = 2802.15392
$ws_ALC.Range("N129").Value = -17342.5
$ws_ALC.Range("H132").Value = 12270.25
$ws_ALC.Range("I132").Value = 13627.4
$ws_ALC.Range("K132").Value = 40882.2
$ws_ALC.Range("M132").Value = -38352.2
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H97").Value = 1322.9412
$ws_ARM.Range("I97").Value = 1186.875
$ws_ARM.Range("J97").Value = 3500
$ws_ARM.Range("K97").Value = 1186.875
$ws_ARM.Range("L97").Value = 3500
$ws_ARM.Range("M97").Value = -690.875
$ws_ARM.Range("N97").Value = -4492
$ws_ARM.Range("H102").Value = 3136.6667
$ws_ARM.Range("I102").Value = 2410
$ws_ARM.Range("J102").Value = 3500
$ws_ARM.Range("K102").Value = 2410
$ws_ARM.Range("L102").Value = 3500
$ws_ARM.Range("M102").Value = -788
$ws_ARM.Range("N102").Value = -6744
$ws_ARM.Range("H132").Value = 2549.578
$ws_ARM.Range("I132").Value = 2171.5642
$ws_ARM.Range("J132").Value = 5006.6665
$ws_ARM.Range("K132").Value = 6514.692599999999
$ws_ARM.Range("L132").Value = 15019.9995
$ws_ARM.Range("M132").Value = -3984.692599999999
$ws_ARM.Range("N132").Value = -20079.9995
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H22").Value = 313
$ws_BSM.Range("I22").Value = 313
$ws_BSM.Range("K22").Value = 313
$ws_BSM.Range("M22").Value = -140
$ws_BSM.Range("H94").Value = 2130.1428
$ws_BSM.Range("I94").Value = 2052.2104
$ws_BSM.Range("J94").Value = 2294.6667
$ws_BSM.Range("K94").Value = 2052.2104
$ws_BSM.Range("L94").Value = 2294.6667
$ws_BSM.Range("M94").Value = -1601.2104
$ws_BSM.Range("N94").Value = -3196.6667
$ws_BSM.Range("H99").Value = 3046.0435
$ws_BSM.Range("I99").Value = 2522.0715
$ws_BSM.Range("J99").Value = 3861.111
$ws_BSM.Range("K99").Value = 2522.0715
$ws_BSM.Range("L99").Value = 3861.111
$ws_BSM.Range("M99").Value = -1024.0715
$ws_BSM.Range("N99").Value = -6857.111
$ws_BSM.Range("H134").Value = 7939107.5
$ws_BSM.Range("I134").Value = 11907161
$ws_BSM.Range("K134").Value = 35721483
$ws_BSM.Range("M134").Value = -35718948
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H133").Value = 200000
$ws_CRP.Range("J133").Value = 200000
$ws_CRP.Range("L133").Value = 200000
$ws_CRP.Range("N133").Value = -205060
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H61").Value = 451.4
$ws_CUL.Range("I61").Value = 200
$ws_CUL.Range("J61").Value = 514.25
$ws_CUL.Range("K61").Value = 600
$ws_CUL.Range("L61").Value = 1542.75
$ws_CUL.Range("M61").Value = -385
$ws_CUL.Range("N61").Value = -1972.75
$ws_CUL.Range("H98").Value = 785.7143
$ws_CUL.Range("I98").Value = 795.5
$ws_CUL.Range("J98").Value = 781.8
$ws_CUL.Range("K98").Value = 2386.5
$ws_CUL.Range("L98").Value = 2345.4
$ws_CUL.Range("M98").Value = -888.5
$ws_CUL.Range("N98").Value = -5341.4
$ws_CUL.Range("H132").Value = 1096
$ws_CUL.Range("J132").Value = 1000
$ws_CUL.Range("L132").Value = 9000
$ws_CUL.Range("N132").Value = -14060
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H132").Value = 2410.5454
$ws_GSM.Range("I132").Value = 1498.2858
$ws_GSM.Range("J132").Value = 4007
$ws_GSM.Range("K132").Value = 4494.857400000001
$ws_GSM.Range("L132").Value = 12021
$ws_GSM.Range("M132").Value = -1964.857400000001
$ws_GSM.Range("N132").Value = -17081
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H22").Value = 3392.4614
$ws_LTW.Range("I22").Value = 3815.6667
$ws_LTW.Range("J22").Value = 3029.7144
$ws_LTW.Range("K22").Value = 3815.6667
$ws_LTW.Range("L22").Value = 3029.7144
$ws_LTW.Range("M22").Value = -3520.6667
$ws_LTW.Range("N22").Value = -3619.7144
$ws_LTW.Range("H27").Value = 3392.4614
$ws_LTW.Range("I27").Value = 3815.6667
$ws_LTW.Range("J27").Value = 3029.7144
$ws_LTW.Range("K27").Value = 3815.6667
$ws_LTW.Range("L27").Value = 3029.7144
$ws_LTW.Range("M27").Value = -3708.6667
$ws_LTW.Range("N27").Value = -3243.7144
$ws_LTW.Range("H108").Value = 74000
$ws_LTW.Range("J108").Value = 74000
$ws_LTW.Range("L108").Value = 74000
$ws_LTW.Range("N108").Value = -81680
$ws_LTW.Range("H132").Value = 3666.6667
$ws_LTW.Range("I132").Value = 3000
$ws_LTW.Range("K132").Value = 9000
$ws_LTW.Range("M132").Value = -6470
$ws_LTW.Range("H136").Value = 4953.4
$ws_LTW.Range("I136").Value = 3324.25
$ws_LTW.Range("K136").Value = 9972.75
$ws_LTW.Range("M136").Value = -7422.75
$ws_LTW.Range("H139").Value = 80000
$ws_LTW.Range("J139").Value = 80000
$ws_LTW.Range("L139").Value = 80000
$ws_LTW.Range("N139").Value = -90280
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H62").Value = 4888.8887
$ws_WVR.Range("J62").Value = 6375
$ws_WVR.Range("L62").Value = 6375
$ws_WVR.Range("N62").Value = -7623
$ws_WVR.Range("H65").Value = 4888.8887
$ws_WVR.Range("J65").Value = 6375
$ws_WVR.Range("L65").Value = 31875
$ws_WVR.Range("N65").Value = -38115
$ws_WVR.Range("H107").Value = 424.89474
$ws_WVR.Range("I107").Value = 403.5
$ws_WVR.Range("J107").Value = 448.66666
$ws_WVR.Range("K107").Value = 1210.5
$ws_WVR.Range("L107").Value = 1345.99998
$ws_WVR.Range("M107").Value = 709.5
$ws_WVR.Range("N107").Value = -5185.999980000001
$ws_WVR.Range("H122").Value = 4301.162
$ws_WVR.Range("I122").Value = 2421.8386
$ws_WVR.Range("K122").Value = 7265.5158
$ws_WVR.Range("M122").Value = -4815.5158
$ws_WVR.Range("H132").Value = 2468.5144
$ws_WVR.Range("I132").Value = 2315.0908
$ws_WVR.Range("K132").Value = 6945.2724
$ws_WVR.Range("M132").Value = -4415.2724
$ws_WVR.Range("H136").Value = 1848.3334
$ws_WVR.Range("I136").Value = 1285.8182
$ws_WVR.Range("J136").Value = 3395.25
$ws_WVR.Range("K136").Value = 3857.4546
$ws_WVR.Range("L136").Value = 10185.75
$ws_WVR.Range("M136").Value = -1307.4546
$ws_WVR.Range("N136").Value = -15285.75
$ws_WVR.Range("H139").Value = 95000
$ws_WVR.Range("J139").Value = 95000
$ws_WVR.Range("L139").Value = 95000
$ws_WVR.Range("N139").Value = -105280

Write-Output "Applied Adamantoise_Profits updates"